# Apply the cryptos-list refresh described in the commit:
#  - Updated Price (D) / Volume(1h) (E) figures for most rows
#  - Swapped the Solana / Cardano rows (9 <-> 10), including Coin name and Link
#
# Columns D/E are stored as plain text in the workbook (not numbers), so for the
# Price cells whose new value still *looks* like a plain number we pre-format the
# cell as Text ('@') before assigning the value - otherwise Excel's COM automation
# would silently re-interpret the string and store it as a numeric cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @(
    'D4',
    'D5',
    'D6',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D14',
    'D15',
    'D16',
    'D17',
    'D19',
    'D20',
    'D21',
    'D22',
    'D24',
    'D26',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D34',
    'D35',
    'D36',
    'D39',
    'D40',
    'D41',
    'D42',
    'D44',
    'D45',
    'D46',
    'D47',
    'D49',
    'D50',
    'D51'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '29.411.46'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.843.47'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = '240.19'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').Value = '0.6312'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').Value = '0.9995'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = '0.07479'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('B9').Value = 'Solana'
$ws.Range('C9').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D9').Value = '25.13'
$ws.Range('E9').Value = '  +3.11%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').Value = '0.2908'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').Value = '0.07743'
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').Value = '1.844.49'
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('D14').Value = '0.6799'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').Value = '0.00001023'
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('D16').Value = '82.14'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('D17').Value = '6.287'
$ws.Range('E17').Value = '  +3.20%  '
$ws.Range('D18').Value = '29.389.41'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').Value = '229.85'
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('D20').Value = '12.34'
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('D21').Value = '0.9995'
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').Value = '7.435'
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').Value = '158.59'
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('D26').Value = '0.1355'
$ws.Range('E26').Value = '  -2.17%  '
$ws.Range('D27').Value = '17.46'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').Value = '0.06609'
$ws.Range('E28').Value = '  +15.94%  '
$ws.Range('D29').Value = '1.431'
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('D30').Value = '1.489'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('D31').Value = '4.077'
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('D32').Value = '4.057'
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('D34').Value = '1.142'
$ws.Range('E34').Value = '  -0.57%  '
$ws.Range('D35').Value = '0.6982'
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('D36').Value = '2.576'
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('E37').Value = '  +2.23%  '
$ws.Range('D38').Value = '1.250.34'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('D39').Value = '2.817'
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').Value = '6.790'
$ws.Range('E40').Value = '  +4.29%  '
$ws.Range('D41').Value = '0.9349'
$ws.Range('E41').Value = '  +3.43%  '
$ws.Range('D42').Value = '0.9995'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').Value = '1.991.09'
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('D44').Value = '101.08'
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('D45').Value = '65.54'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').Value = '7.070'
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').Value = '1.722'
$ws.Range('E47').Value = '  +4.01%  '
$ws.Range('E48').Value = '  +2.88%  '
$ws.Range('D49').Value = '9.072'
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').Value = '0.1149'
$ws.Range('E50').Value = '  -1.09%  '
$ws.Range('D51').Value = '0.3909'
$ws.Range('E51').Value = '  -0.75%  '
